$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '26.159.40'
$ws.Cells.Item(2, 5).Value = '  +3.49%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '1.605.07'
$ws.Cells.Item(3, 5).Value = '  +3.44%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  -0.22%  '

# Row 5
$ws.Cells.Item(5, 4).Value = "'212.89"
$ws.Cells.Item(5, 5).Value = '  +3.04%  '

# Row 6
$ws.Cells.Item(6, 5).Value = '  -0.22%  '

# Row 7
$ws.Cells.Item(7, 5).Value = '  +1.94%  '

# Row 8
$ws.Cells.Item(8, 5).Value = '  +2.53%  '

# Row 9
$ws.Cells.Item(9, 5).Value = '  +1.08%  '

# Row 10
$ws.Cells.Item(10, 4).Value = "'18.04"
$ws.Cells.Item(10, 5).Value = '  +1.64%  '

# Row 11
$ws.Cells.Item(11, 5).Value = '  +5.03%  '

# Row 12
$ws.Cells.Item(12, 4).Value = '1.826.97'
$ws.Cells.Item(12, 5).Value = '  +3.42%  '

# Row 13
$ws.Cells.Item(13, 4).Value = '1.597.84'
$ws.Cells.Item(13, 5).Value = '  +3.09%  '

# Row 14
$ws.Cells.Item(14, 5).Value = '  +0.79%  '

# Row 15
$ws.Cells.Item(15, 5).Value = '  +1.56%  '

# Row 16
$ws.Cells.Item(16, 4).Value = '26.153.63'
$ws.Cells.Item(16, 5).Value = '  +3.57%  '

# Row 17
$ws.Cells.Item(17, 4).Value = "'60.55"
$ws.Cells.Item(17, 5).Value = '  +3.34%  '

# Row 18
$ws.Cells.Item(18, 5).Value = '  +2.30%  '

# Row 20
$ws.Cells.Item(20, 4).Value = "'204.04"
$ws.Cells.Item(20, 5).Value = '  +10.01%  '

# Row 21
$ws.Cells.Item(21, 4).Value = "'4.23"
$ws.Cells.Item(21, 5).Value = '  +3.30%  '

# Row 22
$ws.Cells.Item(22, 4).Value = "'9.30"
$ws.Cells.Item(22, 5).Value = '  +0.54%  '

# Row 23
$ws.Cells.Item(23, 5).Value = '  +2.94%  '

# Row 24
$ws.Cells.Item(24, 5).Value = '  +12.92%  '

# Row 25
$ws.Cells.Item(25, 4).Value = "'141.83"
$ws.Cells.Item(25, 5).Value = '  +2.12%  '

# Row 26
$ws.Cells.Item(26, 5).Value = '  -0.27%  '

# Row 27
$ws.Cells.Item(27, 5).Value = '  -4.23%  '

# Row 28
$ws.Cells.Item(28, 4).Value = "'15.19"
$ws.Cells.Item(28, 5).Value = '  +2.50%  '

# Row 29
$ws.Cells.Item(29, 5).Value = '  +0.85%  '

# Row 30
$ws.Cells.Item(30, 5).Value = '  +1.95%  '

# Row 31
$ws.Cells.Item(31, 5).Value = '  +1.95%  '

# Row 32
$ws.Cells.Item(32, 5).Value = '  +3.17%  '

# Row 33
$ws.Cells.Item(33, 5).Value = '  +1.20%  '

# Row 34
$ws.Cells.Item(34, 5).Value = '  +1.52%  '

# Row 35
$ws.Cells.Item(35, 5).Value = '  +1.69%  '

# Row 36
$ws.Cells.Item(36, 4).Value = "'0.0166"
$ws.Cells.Item(36, 5).Value = '  +11.24%  '

# Row 37
$ws.Cells.Item(37, 4).Value = '1.124.62'
$ws.Cells.Item(37, 5).Value = '  +3.72%  '

# Row 38
$ws.Cells.Item(38, 5).Value = '  +0.11%  '

# Row 39
$ws.Cells.Item(39, 5).Value = '  +3.35%  '

# Row 40
$ws.Cells.Item(40, 5).Value = '  +2.64%  '

# Row 41
$ws.Cells.Item(41, 4).Value = "'0.493"
$ws.Cells.Item(41, 5).Value = '  +0.20%  '

# Row 42
$ws.Cells.Item(42, 5).Value = '  -2.37%  '

# Row 43
$ws.Cells.Item(43, 4).Value = '1.738.65'
$ws.Cells.Item(43, 5).Value = '  +3.41%  '

# Row 44
$ws.Cells.Item(44, 5).Value = '  +1.97%  '

# Row 45
$ws.Cells.Item(45, 5).Value = '  +0.23%  '

# Row 46
$ws.Cells.Item(46, 4).Value = "'1.52"
$ws.Cells.Item(46, 5).Value = '  +4.57%  '

# Row 47
$ws.Cells.Item(47, 4).Value = "'53.54"
$ws.Cells.Item(47, 5).Value = '  +2.54%  '

# Row 48
$ws.Cells.Item(48, 4).Value = "'0.0505"
$ws.Cells.Item(48, 5).Value = '  +0.48%  '

# Row 49
$ws.Cells.Item(49, 4).Value = "'0.409"
$ws.Cells.Item(49, 5).Value = '  +1.23%  '

# Row 50
$ws.Cells.Item(50, 5).Value = '  -0.18%  '

# Row 51
$ws.Cells.Item(51, 2).Value = 'EnergySwap'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(51, 4).Value = "'7.18"
$ws.Cells.Item(51, 5).Value = '  +0.37%  '
